# Updated symbol list on Wed Feb  1 22:17:14 UTC 2023 with GitHub Actions
#
# This script applies the price (column D) and 1h-volume-change (column E)
# refresh recorded in the commit. Values are plain text in the sheet (no
# numeric formatting is applied to D/E), so each assignment is prefixed
# with a leading apostrophe to force a text entry and avoid Excel
# reinterpreting numeric-looking strings (e.g. "317.47") or percent-looking
# strings (e.g. "1.83%") as Number/Percentage values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'317.47"
$ws.Range("E2").Value = "'1.83%"
$ws.Range("D3").Value = "'37.99"
$ws.Range("E3").Value = "'1.73%"
$ws.Range("D4").Value = "'5.177"
$ws.Range("E4").Value = "'0.95%"
$ws.Range("D5").Value = "'0.07984"
$ws.Range("E5").Value = "'2.11%"
$ws.Range("D6").Value = "'8.522"
$ws.Range("E6").Value = "'2.93%"
$ws.Range("D7").Value = "'1.939"
$ws.Range("E7").Value = "'1.84%"
$ws.Range("D8").Value = "'2.986"
$ws.Range("E8").Value = "'5.28%"
$ws.Range("D9").Value = "'0.9418"
$ws.Range("E9").Value = "'2.44%"
$ws.Range("D10").Value = "'0.1243"
$ws.Range("E10").Value = "'3.65%"
$ws.Range("D11").Value = "'0.1938"
$ws.Range("E11").Value = "'1.29%"
$ws.Range("D12").Value = "'0.09049"
$ws.Range("E12").Value = "'0.39%"
$ws.Range("D13").Value = "'0.03408"
$ws.Range("E13").Value = "'1.84%"
$ws.Range("D14").Value = "'0.09553"
$ws.Range("E14").Value = "'-0.47%"
$ws.Range("D15").Value = "'0.001369"
$ws.Range("E15").Value = "'-0.86%"
$ws.Range("D16").Value = "'0.005975"
$ws.Range("E16").Value = "'4.64%"
$ws.Range("E17").Value = "'-2.81%"
$ws.Range("D18").Value = "'4.476"
$ws.Range("E18").Value = "'1.05%"
$ws.Range("E19").Value = "'2.18%"
$ws.Range("D20").Value = "'6.562"
$ws.Range("E20").Value = "'24.59%"
$ws.Range("E21").Value = "'1.84%"
$ws.Range("D22").Value = "'0.2303"
$ws.Range("E22").Value = "'-11.26%"
$ws.Range("D23").Value = "'0.04356"
$ws.Range("E23").Value = "'-0.43%"
$ws.Range("E24").Value = "'-2.09%"
$ws.Range("D25").Value = "'0.004423"
$ws.Range("E25").Value = "'-5.29%"
$ws.Range("D26").Value = "'0.0001324"
$ws.Range("E26").Value = "'-2.91%"
$ws.Range("D27").Value = "'0.0003977"
$ws.Range("E27").Value = "'-0.50%"
$ws.Range("D39").Value = "'0.02429"
$ws.Range("E39").Value = "'7.29%"
$ws.Range("D40").Value = "'0.05163"
$ws.Range("E40").Value = "'2.35%"
$ws.Range("D41").Value = "'0.007452"
$ws.Range("E41").Value = "'-0.29%"
$ws.Range("D42").Value = "'0.1400"
$ws.Range("E42").Value = "'3.77%"
$ws.Range("D43").Value = "'0.008491"
$ws.Range("E43").Value = "'-6.39%"
$ws.Range("D44").Value = "'0.002100"
$ws.Range("E44").Value = "'7.42%"
$ws.Range("D45").Value = "'0.008751"
$ws.Range("E45").Value = "'-5.97%"
$ws.Range("E46").Value = "'-1.20%"
$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("E47").Value = "'-0.57%"
$ws.Range("D48").Value = "'0.002854"
$ws.Range("E48").Value = "'-13.69%"
$ws.Range("D49").Value = "'0.001684"
$ws.Range("E49").Value = "'68.10%"
$ws.Range("D50").Value = "'0.00002093"
$ws.Range("E50").Value = "'-0.57%"
$ws.Range("D51").Value = "'0.0001994"
$ws.Range("E51").Value = "'-0.57%"
